$d = $word.ActiveDocument

# The bold run "DOCX, DOC, PDF, HTML, XPS, R" was split from "TF and TXT"
# by a leftover "_GoBack" bookmark. Replacing the full (already
# concatenated) phrase merges the two runs back into one and drops the
# now-unneeded bookmarkStart/bookmarkEnd pair.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("DOCX, DOC, PDF, HTML, XPS, RTF and TXT", $true, $false, $false, $false, $false, $true, 1, $false, `
               "DOCX, DOC, PDF, HTML, XPS, RTF and TXT", 2) | Out-Null
